$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows at row 972 (shifts existing rows 972-1044 down to 981-1053)
$ws.Rows("972:980").Insert()

# Populate the newly inserted rows with historical data (2019-11-18 through 2019-11-28)
$ws.Cells.Item(972,1).Value = 1574035200
$ws.Cells.Item(972,2).NumberFormat = "@"
$ws.Cells.Item(972,2).Value = "2019-11-18"
$ws.Cells.Item(972,2).Style = "Normal"
$ws.Cells.Item(972,3).NumberFormat = "@"
$ws.Cells.Item(972,3).Value = "0180"
$ws.Cells.Item(972,3).Style = "Normal"
$ws.Cells.Item(972,4).Value = "KTC"
$ws.Cells.Item(972,5).Value = 0.15
$ws.Cells.Item(972,6).Value = 0.155
$ws.Cells.Item(972,7).Value = 0.15
$ws.Cells.Item(972,8).Value = 0.155
$ws.Cells.Item(972,9).Value = 83000

$ws.Cells.Item(973,1).Value = 1574121600
$ws.Cells.Item(973,2).NumberFormat = "@"
$ws.Cells.Item(973,2).Value = "2019-11-19"
$ws.Cells.Item(973,2).Style = "Normal"
$ws.Cells.Item(973,3).NumberFormat = "@"
$ws.Cells.Item(973,3).Value = "0180"
$ws.Cells.Item(973,3).Style = "Normal"
$ws.Cells.Item(973,4).Value = "KTC"
$ws.Cells.Item(973,5).Value = 0.155
$ws.Cells.Item(973,6).Value = 0.155
$ws.Cells.Item(973,7).Value = 0.155
$ws.Cells.Item(973,8).Value = 0.155
$ws.Cells.Item(973,9).Value = 439700

$ws.Cells.Item(974,1).Value = 1574208000
$ws.Cells.Item(974,2).NumberFormat = "@"
$ws.Cells.Item(974,2).Value = "2019-11-20"
$ws.Cells.Item(974,2).Style = "Normal"
$ws.Cells.Item(974,3).NumberFormat = "@"
$ws.Cells.Item(974,3).Value = "0180"
$ws.Cells.Item(974,3).Style = "Normal"
$ws.Cells.Item(974,4).Value = "KTC"
$ws.Cells.Item(974,5).Value = 0.155
$ws.Cells.Item(974,6).Value = 0.155
$ws.Cells.Item(974,7).Value = 0.155
$ws.Cells.Item(974,8).Value = 0.155
$ws.Cells.Item(974,9).Value = 50000

$ws.Cells.Item(975,1).Value = 1574294400
$ws.Cells.Item(975,2).NumberFormat = "@"
$ws.Cells.Item(975,2).Value = "2019-11-21"
$ws.Cells.Item(975,2).Style = "Normal"
$ws.Cells.Item(975,3).NumberFormat = "@"
$ws.Cells.Item(975,3).Value = "0180"
$ws.Cells.Item(975,3).Style = "Normal"
$ws.Cells.Item(975,4).Value = "KTC"
$ws.Cells.Item(975,5).Value = 0.155
$ws.Cells.Item(975,6).Value = 0.165
$ws.Cells.Item(975,7).Value = 0.155
$ws.Cells.Item(975,8).Value = 0.165
$ws.Cells.Item(975,9).Value = 1165600

$ws.Cells.Item(976,1).Value = 1574380800
$ws.Cells.Item(976,2).NumberFormat = "@"
$ws.Cells.Item(976,2).Value = "2019-11-22"
$ws.Cells.Item(976,2).Style = "Normal"
$ws.Cells.Item(976,3).NumberFormat = "@"
$ws.Cells.Item(976,3).Value = "0180"
$ws.Cells.Item(976,3).Style = "Normal"
$ws.Cells.Item(976,4).Value = "KTC"
$ws.Cells.Item(976,5).Value = 0.165
$ws.Cells.Item(976,6).Value = 0.165
$ws.Cells.Item(976,7).Value = 0.16
$ws.Cells.Item(976,8).Value = 0.16
$ws.Cells.Item(976,9).Value = 267000

$ws.Cells.Item(977,1).Value = 1574640000
$ws.Cells.Item(977,2).NumberFormat = "@"
$ws.Cells.Item(977,2).Value = "2019-11-25"
$ws.Cells.Item(977,2).Style = "Normal"
$ws.Cells.Item(977,3).NumberFormat = "@"
$ws.Cells.Item(977,3).Value = "0180"
$ws.Cells.Item(977,3).Style = "Normal"
$ws.Cells.Item(977,4).Value = "KTC"
$ws.Cells.Item(977,5).Value = 0.165
$ws.Cells.Item(977,6).Value = 0.165
$ws.Cells.Item(977,7).Value = 0.16
$ws.Cells.Item(977,8).Value = 0.16
$ws.Cells.Item(977,9).Value = 233800

$ws.Cells.Item(978,1).Value = 1574726400
$ws.Cells.Item(978,2).NumberFormat = "@"
$ws.Cells.Item(978,2).Value = "2019-11-26"
$ws.Cells.Item(978,2).Style = "Normal"
$ws.Cells.Item(978,3).NumberFormat = "@"
$ws.Cells.Item(978,3).Value = "0180"
$ws.Cells.Item(978,3).Style = "Normal"
$ws.Cells.Item(978,4).Value = "KTC"
$ws.Cells.Item(978,5).Value = 0.155
$ws.Cells.Item(978,6).Value = 0.155
$ws.Cells.Item(978,7).Value = 0.155
$ws.Cells.Item(978,8).Value = 0.155
$ws.Cells.Item(978,9).Value = 40000

$ws.Cells.Item(979,1).Value = 1574812800
$ws.Cells.Item(979,2).NumberFormat = "@"
$ws.Cells.Item(979,2).Value = "2019-11-27"
$ws.Cells.Item(979,2).Style = "Normal"
$ws.Cells.Item(979,3).NumberFormat = "@"
$ws.Cells.Item(979,3).Value = "0180"
$ws.Cells.Item(979,3).Style = "Normal"
$ws.Cells.Item(979,4).Value = "KTC"
$ws.Cells.Item(979,5).Value = 0.16
$ws.Cells.Item(979,6).Value = 0.165
$ws.Cells.Item(979,7).Value = 0.16
$ws.Cells.Item(979,8).Value = 0.165
$ws.Cells.Item(979,9).Value = 233000

$ws.Cells.Item(980,1).Value = 1574899200
$ws.Cells.Item(980,2).NumberFormat = "@"
$ws.Cells.Item(980,2).Value = "2019-11-28"
$ws.Cells.Item(980,2).Style = "Normal"
$ws.Cells.Item(980,3).NumberFormat = "@"
$ws.Cells.Item(980,3).Value = "0180"
$ws.Cells.Item(980,3).Style = "Normal"
$ws.Cells.Item(980,4).Value = "KTC"
$ws.Cells.Item(980,5).Value = 0.16
$ws.Cells.Item(980,6).Value = 0.16
$ws.Cells.Item(980,7).Value = 0.16
$ws.Cells.Item(980,8).Value = 0.16
$ws.Cells.Item(980,9).Value = 253100

